$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.326.35"
$ws.Range("E2").Value = "  -2.27%  "
$ws.Range("D3").Value = "1.933.90"
$ws.Range("E3").Value = "  -2.21%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.69%  "
$ws.Range("D5").Value = "'250.93"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("D6").Value = "'0.7124"
$ws.Range("E6").Value = "  -2.37%  "
$ws.Range("E7").Value = "  -0.52%  "
$ws.Range("D8").Value = "'0.3304"
$ws.Range("E8").Value = "  -1.87%  "
$ws.Range("D9").Value = "'27.73"
$ws.Range("E9").Value = "  +1.11%  "
$ws.Range("D10").Value = "'0.07289"
$ws.Range("E10").Value = "  +2.64%  "
$ws.Range("D11").Value = "'0.8056"
$ws.Range("E11").Value = "  -2.26%  "
$ws.Range("D12").Value = "'0.08095"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").Value = "1.932.98"
$ws.Range("E13").Value = "  -2.34%  "
$ws.Range("E14").Value = "  -1.54%  "
$ws.Range("D15").Value = "'94.62"
$ws.Range("E15").Value = "  -4.14%  "
$ws.Range("E16").Value = "  -1.04%  "
$ws.Range("D17").Value = "30.319.85"
$ws.Range("E17").Value = "  -2.37%  "
$ws.Range("D18").Value = "'253.34"
$ws.Range("E18").Value = "  -5.06%  "
$ws.Range("D19").Value = "'0.000008190"
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").Value = "'5.789"
$ws.Range("E20").Value = "  -4.64%  "
$ws.Range("D21").Value = "2.187.99"
$ws.Range("E21").Value = "  -2.64%  "
$ws.Range("E22").Value = "  -0.55%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  -0.79%  "
$ws.Range("D24").Value = "'6.975"
$ws.Range("E24").Value = "  -1.15%  "
$ws.Range("D25").Value = "'9.759"
$ws.Range("E25").Value = "  -1.72%  "
$ws.Range("D26").Value = "'165.16"
$ws.Range("E26").Value = "  +2.04%  "
$ws.Range("D27").Value = "'2.360"
$ws.Range("E27").Value = "  +0.54%  "
$ws.Range("D28").Value = "'19.32"
$ws.Range("E28").Value = "  -1.63%  "
$ws.Range("D29").Value = "'0.1288"
$ws.Range("E29").Value = "  -2.23%  "
$ws.Range("D30").Value = "'1.350"
$ws.Range("E30").Value = "  -1.70%  "
$ws.Range("D31").Value = "'1.543"
$ws.Range("E31").Value = "  -2.97%  "
$ws.Range("D32").Value = "'4.419"
$ws.Range("E32").Value = "  -4.21%  "
$ws.Range("D33").Value = "'4.178"
$ws.Range("E33").Value = "  -5.04%  "
$ws.Range("D34").Value = "'0.05197"
$ws.Range("E34").Value = "  -1.51%  "
$ws.Range("D35").Value = "'1.264"
$ws.Range("E35").Value = "  -0.48%  "
$ws.Range("D36").Value = "'0.7461"
$ws.Range("E36").Value = "  -3.85%  "
$ws.Range("D37").Value = "'2.785"
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("D38").Value = "'0.01967"
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("D39").Value = "'2.804"
$ws.Range("E39").Value = "  -2.70%  "
$ws.Range("D40").Value = "'78.89"
$ws.Range("E40").Value = "  -5.12%  "
$ws.Range("D41").Value = "'6.423"
$ws.Range("E41").Value = "  -4.79%  "
$ws.Range("D42").Value = "'0.4527"
$ws.Range("E42").Value = "  -1.81%  "
$ws.Range("D43").Value = "'2.019"
$ws.Range("E43").Value = "  -3.21%  "
$ws.Range("D44").Value = "'0.8443"
$ws.Range("E44").Value = "  -0.83%  "
$ws.Range("D45").Value = "'1.000"
$ws.Range("E45").Value = "  -0.55%  "
$ws.Range("D46").Value = "'101.61"
$ws.Range("E46").Value = "  -2.75%  "
$ws.Range("D47").Value = "'9.762"
$ws.Range("E47").Value = "  -2.80%  "
$ws.Range("D48").Value = "'7.455"
$ws.Range("E48").Value = "  -2.15%  "
$ws.Range("D49").Value = "'36.75"
$ws.Range("E49").Value = "  -0.99%  "
$ws.Range("D50").Value = "'0.4172"
$ws.Range("E50").Value = "  -2.60%  "
$ws.Range("D51").Value = "'0.06033"
$ws.Range("E51").Value = "  +0.04%  "
